$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.490.58'
$ws.Range('E2').Value = '  -0.04%  '

# Row 3
$ws.Range('D3').Value = '3.897.81'
$ws.Range('E3').Value = '  -0.03%  '

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$ws.Range('D5').Value = '''602.88'
$ws.Range('E5').Value = '  +0.12%  '

# Row 6
$ws.Range('D6').Value = '''169.13'
$ws.Range('E6').Value = '  +1.80%  '

# Row 7
$ws.Range('D7').Value = '3.897.10'
$ws.Range('E7').Value = '  +0.06%  '

# Row 8
$ws.Range('E8').Value = '  +0.05%  '

# Row 9
$ws.Range('E9').Value = '  +0.17%  '

# Row 10
$ws.Range('E10').Value = '  +0.10%  '

# Row 11
$ws.Range('E11').Value = '  +0.43%  '

# Row 12
$ws.Range('E12').Value = '  -0.25%  '

# Row 13
$ws.Range('D13').Value = '''0.0000254'
$ws.Range('E13').Value = '  -0.18%  '

# Row 14
$ws.Range('D14').Value = '''37.09'
$ws.Range('E14').Value = '  -0.64%  '

# Row 15
$ws.Range('D15').Value = '4.551.31'
$ws.Range('E15').Value = '  +0.21%  '

# Row 16
$ws.Range('D16').Value = '3.900.80'
$ws.Range('E16').Value = '  +0.17%  '

# Row 17
$ws.Range('D17').Value = '68.445.48'
$ws.Range('E17').Value = '  -0.18%  '

# Row 18
$ws.Range('D18').Value = '''18.13'
$ws.Range('E18').Value = '  +5.22%  '

# Row 19
$ws.Range('D19').Value = '''7.42'
$ws.Range('E19').Value = '  -0.83%  '

# Row 20
$ws.Range('E20').Value = '  +0.23%  '

# Row 21
$ws.Range('D21').Value = '''10.86'
$ws.Range('E21').Value = '  -1.59%  '

# Row 22
$ws.Range('D22').Value = '''473.13'
$ws.Range('E22').Value = '  -2.98%  '

# Row 23
$ws.Range('D23').Value = '''0.739'
$ws.Range('E23').Value = '  +2.11%  '

# Row 24
$ws.Range('D24').Value = '''0.0000169'
$ws.Range('E24').Value = '  +1.58%  '

# Row 25
$ws.Range('D25').Value = '''83.78'
$ws.Range('E25').Value = '  -0.97%  '

# Row 26
$ws.Range('E26').Value = '  +1.08%  '

# Row 27
$ws.Range('D27').Value = '''12.24'
$ws.Range('E27').Value = '  +1.39%  '

# Row 28
$ws.Range('E28').Value = '  +0.10%  '

# Row 29
$ws.Range('D29').Value = '''10.01'
$ws.Range('E29').Value = '  -0.94%  '

# Row 30
$ws.Range('E30').Value = '  +1.25%  '

# Row 31
$ws.Range('D31').Value = '4.047.55'
$ws.Range('E31').Value = '  +0.08%  '

# Row 32
$ws.Range('E32').Value = '  +1.73%  '

# Row 33
$ws.Range('D33').Value = '''31.51'
$ws.Range('E33').Value = '  -1.13%  '

# Row 34
$ws.Range('D34').Value = '''2.31'
$ws.Range('E34').Value = '  -2.71%  '

# Row 35
$ws.Range('D35').Value = '''9.42'
$ws.Range('E35').Value = '  +1.31%  '

# Row 36
$ws.Range('D36').Value = '3.869.87'
$ws.Range('E36').Value = '  +0.48%  '

# Row 37
$ws.Range('E37').Value = '  -1.76%  '

# Row 38
$ws.Range('D38').Value = '''3.68'
$ws.Range('E38').Value = '  +15.84%  '

# Row 39
$ws.Range('D39').Value = '''1.04'
$ws.Range('E39').Value = '  -0.15%  '

# Row 40
$ws.Range('E40').Value = '  +2.00%  '

# Row 41
$ws.Range('E41').Value = '  -0.30%  '

# Row 42
$ws.Range('E42').Value = '  +0.13%  '

# Row 43
$ws.Range('E43').Value = '  -0.42%  '

# Row 44
$ws.Range('B44').Value = 'FLOKI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D44').Value = '''0.000304'
$ws.Range('E44').Value = '  +15.81%  '

# Row 45
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '''429.03'
$ws.Range('E45').Value = '  +0.03%  '

# Row 46
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '''2.00'
$ws.Range('E46').Value = '  +0.91%  '

# Row 47
$ws.Range('E47').Value = '  +0.02%  '

# Row 48
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '''47.25'
$ws.Range('E48').Value = '  -2.17%  '

# Row 49
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '''8.61'
$ws.Range('E49').Value = '  +1.09%  '

# Row 50
$ws.Range('D50').Value = '''26.98'
$ws.Range('E50').Value = '  +5.09%  '

# Row 51
$ws.Range('D51').Value = '''143.77'
$ws.Range('E51').Value = '  +0.97%  '

